$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 442190.3
$ws.Range("J17").Value = 442190.3
$ws.Range("L17").Value = 1326570.9
$ws.Range("N17").Value = -1326906.9
$ws.Range("H51").Value = 11299.167
$ws.Range("I51").Value = 3999.3333
$ws.Range("J51").Value = 13732.444
$ws.Range("K51").Value = 3999.3333
$ws.Range("L51").Value = 13732.444
$ws.Range("M51").Value = -3515.3333
$ws.Range("N51").Value = -14700.444
$ws.Range("H52").Value = 1224.5
$ws.Range("I52").Value = 950
$ws.Range("J52").Value = 1499
$ws.Range("K52").Value = 2850
$ws.Range("L52").Value = 4497
$ws.Range("M52").Value = -2690
$ws.Range("N52").Value = -4817
$ws.Range("H94").Value = 7025.75
$ws.Range("I94").Value = 4926.5
$ws.Range("J94").Value = 9125
$ws.Range("K94").Value = 4926.5
$ws.Range("L94").Value = 9125
$ws.Range("M94").Value = -4475.5
$ws.Range("N94").Value = -10027
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 12500
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 12500
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 12500
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -12788
$ws.Range("H63").Value = 4020
$ws.Range("I63").Value = 2900
$ws.Range("J63").Value = 8500
$ws.Range("K63").Value = 2900
$ws.Range("L63").Value = 8500
$ws.Range("M63").Value = -2214
$ws.Range("N63").Value = -9872
$ws.Range("H66").Value = 4020
$ws.Range("I66").Value = 2900
$ws.Range("J66").Value = 8500
$ws.Range("K66").Value = 14500
$ws.Range("L66").Value = 42500
$ws.Range("M66").Value = -11068
$ws.Range("N66").Value = -49364
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 5756.0625
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 5756.0625
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 5756.0625
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -6210.0625
$ws.Range("H19").Value = 12990.091
$ws.Range("J19").Value = 12990.091
$ws.Range("L19").Value = 12990.091
$ws.Range("N19").Value = -13336.091
$ws.Range("H35").Value = 15722.223
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 15722.223
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 15722.223
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -16342.223
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 28151.2
$ws.Range("I10").Value = 246.66667
$ws.Range("K10").Value = 246.66667
$ws.Range("M10").Value = -107.66667
$ws.Range("H99").Value = 4238.5
$ws.Range("I99").Value = 4352.375
$ws.Range("J99").Value = 4124.625
$ws.Range("K99").Value = 4352.375
$ws.Range("L99").Value = 4124.625
$ws.Range("M99").Value = -2854.375
$ws.Range("N99").Value = -7120.625
$ws.Range("H126").Value = 4238.5
$ws.Range("I126").Value = 4352.375
$ws.Range("J126").Value = 4124.625
$ws.Range("K126").Value = 13057.125
$ws.Range("L126").Value = 12373.875
$ws.Range("M126").Value = -10587.125
$ws.Range("N126").Value = -17313.875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2861.2727
$ws.Range("J3").Value = 4026.2856
$ws.Range("L3").Value = 12078.8568
$ws.Range("N3").Value = -12302.8568
$ws.Range("H5").Value = 7462.3335
$ws.Range("I5").Value = 899
$ws.Range("J5").Value = 9849
$ws.Range("K5").Value = 2697
$ws.Range("L5").Value = 29547
$ws.Range("M5").Value = -2585
$ws.Range("N5").Value = -29771
$ws.Range("H38").Value = 111.36364
$ws.Range("I38").Value = 103.125
$ws.Range("J38").Value = 133.33333
$ws.Range("K38").Value = 309.375
$ws.Range("L38").Value = 399.99999
$ws.Range("M38").Value = 37.625
$ws.Range("N38").Value = -1093.99999
$ws.Range("H113").Value = 561.9048
$ws.Range("I113").Value = 430
$ws.Range("J113").Value = 603.125
$ws.Range("K113").Value = 1290
$ws.Range("L113").Value = 1809.375
$ws.Range("M113").Value = 880
$ws.Range("N113").Value = -6149.375
$ws.Range("H122").Value = 339.2
$ws.Range("I122").Value = 294.83334
$ws.Range("J122").Value = 516.6667
$ws.Range("K122").Value = 2653.50006
$ws.Range("L122").Value = 4650.0003
$ws.Range("M122").Value = -203.5000600000003
$ws.Range("N122").Value = -9550.0003
$ws.Range("H135").Value = 7462.3335
$ws.Range("I135").Value = 899
$ws.Range("J135").Value = 9849
$ws.Range("K135").Value = 8091
$ws.Range("L135").Value = 88641
$ws.Range("M135").Value = -5556
$ws.Range("N135").Value = -93711
$ws.Range("H140").Value = 13903936
$ws.Range("I140").Value = 3276
$ws.Range("J140").Value = 37071704
$ws.Range("K140").Value = 9828
$ws.Range("L140").Value = 111215112
$ws.Range("M140").Value = -4648
$ws.Range("N140").Value = -111225472
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4036.6511
$ws.Range("I70").Value = 3954.4243
$ws.Range("J70").Value = 4308
$ws.Range("K70").Value = 3954.4243
$ws.Range("L70").Value = 4308
$ws.Range("M70").Value = -3684.4243
$ws.Range("N70").Value = -4848
$ws.Range("H73").Value = 4036.6511
$ws.Range("I73").Value = 3954.4243
$ws.Range("J73").Value = 4308
$ws.Range("K73").Value = 3954.4243
$ws.Range("L73").Value = 4308
$ws.Range("M73").Value = -3018.4243
$ws.Range("N73").Value = -6180
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3070.8333
$ws.Range("I7").Value = 3697
$ws.Range("J7").Value = 2355.2144
$ws.Range("K7").Value = 3697
$ws.Range("L7").Value = 2355.2144
$ws.Range("M7").Value = -3585
$ws.Range("N7").Value = -2579.2144
$ws.Range("H18").Value = 20000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H20").Value = 40953
$ws.Range("J20").Value = 40953
$ws.Range("L20").Value = 40953
$ws.Range("N20").Value = -41405
$ws.Range("H40").Value = 33930.22
$ws.Range("I40").Value = 1975.1875
$ws.Range("J40").Value = 65885.25
$ws.Range("K40").Value = 1975.1875
$ws.Range("L40").Value = 65885.25
$ws.Range("M40").Value = -1839.1875
$ws.Range("N40").Value = -66157.25
$ws.Range("H68").Value = 2163.077
$ws.Range("J68").Value = 2744
$ws.Range("L68").Value = 2744
$ws.Range("N68").Value = -4242
$ws.Range("H71").Value = 2163.077
$ws.Range("J71").Value = 2744
$ws.Range("L71").Value = 13720
$ws.Range("N71").Value = -21208
$ws.Range("H122").Value = 1894.3
$ws.Range("I122").Value = 2149.3333
$ws.Range("J122").Value = 1785
$ws.Range("K122").Value = 6447.999899999999
$ws.Range("L122").Value = 5355
$ws.Range("M122").Value = -3997.999899999999
$ws.Range("N122").Value = -10255
$ws.Range("H126").Value = 3070.8333
$ws.Range("I126").Value = 3697
$ws.Range("J126").Value = 2355.2144
$ws.Range("K126").Value = 11091
$ws.Range("L126").Value = 7065.6432
$ws.Range("M126").Value = -8621
$ws.Range("N126").Value = -12005.6432
$ws.Range("H132").Value = 457593.72
$ws.Range("I132").Value = 2893.3333
$ws.Range("J132").Value = 1431951.8
$ws.Range("K132").Value = 8679.999899999999
$ws.Range("L132").Value = 4295855.4
$ws.Range("M132").Value = -6149.999899999999
$ws.Range("N132").Value = -4300915.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1219.1177
$ws.Range("I126").Value = 1062.5
$ws.Range("J126").Value = 1442.8572
$ws.Range("K126").Value = 3187.5
$ws.Range("L126").Value = 4328.571599999999
$ws.Range("M126").Value = -717.5
$ws.Range("N126").Value = -9268.571599999999
$ws.Range("H132").Value = 1655
$ws.Range("I132").Value = 1342.0454
$ws.Range("K132").Value = 4026.1362
$ws.Range("M132").Value = -1496.1362
